$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.ClearFormats()
}

Set-TextValue $ws.Range("D2") "331.14"
Set-TextValue $ws.Range("E2") "0.24%"
Set-TextValue $ws.Range("D3") "41.33"
Set-TextValue $ws.Range("E3") "0.34%"
Set-TextValue $ws.Range("D4") "5.710"
Set-TextValue $ws.Range("E4") "-0.14%"
Set-TextValue $ws.Range("E5") "4.06%"
Set-TextValue $ws.Range("D6") "8.821"
Set-TextValue $ws.Range("E6") "0.93%"
Set-TextValue $ws.Range("D7") "4.505"
Set-TextValue $ws.Range("E7") "-0.07%"
Set-TextValue $ws.Range("D8") "1.986"
Set-TextValue $ws.Range("E8") "-2.80%"
Set-TextValue $ws.Range("D10") "0.9288"
Set-TextValue $ws.Range("E10") "0.59%"
Set-TextValue $ws.Range("D11") "0.1253"
Set-TextValue $ws.Range("E11") "0.18%"
Set-TextValue $ws.Range("D12") "0.1966"
Set-TextValue $ws.Range("E12") "1.23%"
Set-TextValue $ws.Range("D13") "0.09402"
Set-TextValue $ws.Range("E13") "1.30%"
Set-TextValue $ws.Range("D14") "0.03961"
Set-TextValue $ws.Range("E14") "8.16%"
Set-TextValue $ws.Range("D15") "0.1064"
Set-TextValue $ws.Range("E15") "0.87%"
Set-TextValue $ws.Range("D16") "0.001313"
Set-TextValue $ws.Range("E16") "1.01%"
Set-TextValue $ws.Range("D17") "0.006117"
Set-TextValue $ws.Range("E17") "-0.12%"
Set-TextValue $ws.Range("D18") "3.434"
Set-TextValue $ws.Range("E18") "1.51%"
Set-TextValue $ws.Range("E19") "0.83%"
Set-TextValue $ws.Range("D20") "9.155"
Set-TextValue $ws.Range("D21") "0.1365"
Set-TextValue $ws.Range("E21") "-3.57%"
Set-TextValue $ws.Range("D22") "0.2513"
Set-TextValue $ws.Range("E22") "-5.21%"
Set-TextValue $ws.Range("D23") "0.04409"
Set-TextValue $ws.Range("E23") "-0.43%"
Set-TextValue $ws.Range("D24") "0.001246"
Set-TextValue $ws.Range("E24") "-1.09%"
Set-TextValue $ws.Range("D25") "0.004392"
Set-TextValue $ws.Range("E25") "1.02%"
Set-TextValue $ws.Range("D26") "0.0001192"
Set-TextValue $ws.Range("E26") "-3.97%"
Set-TextValue $ws.Range("D27") "0.0003998"
Set-TextValue $ws.Range("E27") "0.14%"
Set-TextValue $ws.Range("D39") "0.02834"
Set-TextValue $ws.Range("E39") "0.74%"
Set-TextValue $ws.Range("D40") "0.05515"
Set-TextValue $ws.Range("E40") "0.81%"
Set-TextValue $ws.Range("D41") "0.007920"
Set-TextValue $ws.Range("E41") "4.18%"
Set-TextValue $ws.Range("E42") "1.23%"
Set-TextValue $ws.Range("D43") "0.008974"
Set-TextValue $ws.Range("E43") "-9.79%"
Set-TextValue $ws.Range("D44") "0.002083"
Set-TextValue $ws.Range("E44") "-1.74%"
Set-TextValue $ws.Range("D45") "0.01021"
Set-TextValue $ws.Range("E45") "-13.23%"
Set-TextValue $ws.Range("D46") "0.00007315"
Set-TextValue $ws.Range("E46") "8.65%"
Set-TextValue $ws.Range("D47") "0.00000000752"
Set-TextValue $ws.Range("E47") "0.24%"
Set-TextValue $ws.Range("D48") "0.003212"
Set-TextValue $ws.Range("E48") "8.73%"
Set-TextValue $ws.Range("D49") "0.002283"
Set-TextValue $ws.Range("E49") "0.23%"
Set-TextValue $ws.Range("D50") "0.00002105"
Set-TextValue $ws.Range("E50") "0.24%"
Set-TextValue $ws.Range("D51") "0.0002005"
Set-TextValue $ws.Range("E51") "0.24%"
